# Auto-generated Excel COM-interop script
# Applies updated market-price snapshot values (columns H:N) produced by the
# scheduled pricing runner, matching the authoritative diff for this commit.

$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 390.8
$ws.Range("I2").Value = 318
$ws.Range("K2").Value = 318
$ws.Range("M2").Value = -205
$ws.Range("H5").Value = 155.92308
$ws.Range("I5").Value = 206.75
$ws.Range("J5").Value = 74.59999999999999
$ws.Range("K5").Value = 206.75
$ws.Range("L5").Value = 74.59999999999999
$ws.Range("M5").Value = -91.75
$ws.Range("N5").Value = -304.6
$ws.Range("H18").Value = 8242.333000000001
$ws.Range("I18").Value = 8091
$ws.Range("K18").Value = 8091
$ws.Range("M18").Value = -7807
$ws.Range("H112").Value = 4036.0715
$ws.Range("J112").Value = 4065.7693
$ws.Range("L112").Value = 12197.3079
$ws.Range("N112").Value = -14413.3079

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2206.875
$ws.Range("I2").Value = 2044.1111
$ws.Range("J2").Value = 2416.1428
$ws.Range("K2").Value = 2044.1111
$ws.Range("L2").Value = 2416.1428
$ws.Range("M2").Value = -1931.1111
$ws.Range("N2").Value = -2642.1428
$ws.Range("H116").Value = 2206.875
$ws.Range("I116").Value = 2044.1111
$ws.Range("J116").Value = 2416.1428
$ws.Range("K116").Value = 2044.1111
$ws.Range("L116").Value = 2416.1428
$ws.Range("M116").Value = 249.8888999999999
$ws.Range("N116").Value = -7004.1428
$ws.Range("H132").Value = 2792.8125
$ws.Range("I132").Value = 2616.6072
$ws.Range("K132").Value = 7849.821599999999
$ws.Range("M132").Value = -5319.821599999999

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2206.875
$ws.Range("I3").Value = 2044.1111
$ws.Range("J3").Value = 2416.1428
$ws.Range("K3").Value = 2044.1111
$ws.Range("L3").Value = 2416.1428
$ws.Range("M3").Value = -1930.1111
$ws.Range("N3").Value = -2644.1428
$ws.Range("H94").Value = 946.5294
$ws.Range("I94").Value = 1007.0909
$ws.Range("J94").Value = 835.5
$ws.Range("K94").Value = 1007.0909
$ws.Range("L94").Value = 835.5
$ws.Range("M94").Value = -556.0909
$ws.Range("N94").Value = -1737.5
$ws.Range("H105").Value = 5657.4644
$ws.Range("I105").Value = 5612.7334
$ws.Range("J105").Value = 5709.077
$ws.Range("K105").Value = 5612.7334
$ws.Range("L105").Value = 5709.077
$ws.Range("M105").Value = -3865.7334
$ws.Range("N105").Value = -9203.077000000001

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2695.3333
$ws.Range("I31").Value = 1335.6666
$ws.Range("J31").Value = 9493.666999999999
$ws.Range("K31").Value = 1335.6666
$ws.Range("L31").Value = 9493.666999999999
$ws.Range("M31").Value = -1040.6666
$ws.Range("N31").Value = -10083.667
$ws.Range("H34").Value = 2695.3333
$ws.Range("I34").Value = 1335.6666
$ws.Range("J34").Value = 9493.666999999999
$ws.Range("K34").Value = 1335.6666
$ws.Range("L34").Value = 9493.666999999999
$ws.Range("M34").Value = -1133.6666
$ws.Range("N34").Value = -9897.666999999999
$ws.Range("H41").Value = 26229.125
$ws.Range("I41").Value = 13855.6
$ws.Range("K41").Value = 13855.6
$ws.Range("M41").Value = -13427.6
$ws.Range("H50").Value = 50166.168
$ws.Range("I50").Value = 40500
$ws.Range("J50").Value = 54999.25
$ws.Range("K50").Value = 40500
$ws.Range("L50").Value = 54999.25
$ws.Range("M50").Value = -39875
$ws.Range("N50").Value = -56249.25
$ws.Range("H62").Value = 3163.9092
$ws.Range("I62").Value = 2609.5
$ws.Range("J62").Value = 3829.2
$ws.Range("K62").Value = 2609.5
$ws.Range("L62").Value = 3829.2
$ws.Range("M62").Value = -1985.5
$ws.Range("N62").Value = -5077.2
$ws.Range("H65").Value = 3163.9092
$ws.Range("I65").Value = 2609.5
$ws.Range("J65").Value = 3829.2
$ws.Range("K65").Value = 13047.5
$ws.Range("L65").Value = 19146
$ws.Range("M65").Value = -9927.5
$ws.Range("N65").Value = -25386

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 8399.5
$ws.Range("J62").Value = 7666
$ws.Range("L62").Value = 22998
$ws.Range("N62").Value = -24370
$ws.Range("H65").Value = 8399.5
$ws.Range("J65").Value = 7666
$ws.Range("L65").Value = 68994
$ws.Range("N65").Value = -75858
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H126").Value = 2782.25
$ws.Range("I126").Value = 2782.25
$ws.Range("K126").Value = 8346.75
$ws.Range("M126").Value = -3406.75
$ws.Range("H138").Value = 7916
$ws.Range("I138").Value = 7916
$ws.Range("K138").Value = 23748
$ws.Range("M138").Value = -18608

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2963.8
$ws.Range("I80").Value = 2718
$ws.Range("J80").Value = 3332.5
$ws.Range("K80").Value = 2718
$ws.Range("L80").Value = 3332.5
$ws.Range("M80").Value = -1720
$ws.Range("N80").Value = -5328.5
$ws.Range("H83").Value = 2963.8
$ws.Range("I83").Value = 2718
$ws.Range("J83").Value = 3332.5
$ws.Range("K83").Value = 13590
$ws.Range("L83").Value = 16662.5
$ws.Range("M83").Value = -8598
$ws.Range("N83").Value = -26646.5
$ws.Range("H122").Value = 2217.0476
$ws.Range("I122").Value = 2083.3333
$ws.Range("K122").Value = 6249.999899999999
$ws.Range("M122").Value = -3799.999899999999
$ws.Range("H126").Value = 4259.3335
$ws.Range("I126").Value = 2777
$ws.Range("K126").Value = 8331
$ws.Range("M126").Value = -5861

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 323.0909
$ws.Range("I16").Value = 323.0909
$ws.Range("K16").Value = 323.0909
$ws.Range("M16").Value = -153.0909
$ws.Range("H46").Value = 2407.1667
$ws.Range("J46").Value = 3284.3333
$ws.Range("L46").Value = 3284.3333
$ws.Range("N46").Value = -3660.3333
$ws.Range("H122").Value = 2870.7334
$ws.Range("J122").Value = 3500
$ws.Range("L122").Value = 10500
$ws.Range("N122").Value = -15400
$ws.Range("H136").Value = 23812718
$ws.Range("I136").Value = 2922.7693
$ws.Range("K136").Value = 8768.3079
$ws.Range("M136").Value = -6218.3079

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 434.75
$ws.Range("I113").Value = 449.66666
$ws.Range("J113").Value = 390
$ws.Range("K113").Value = 1348.99998
$ws.Range("L113").Value = 1170
$ws.Range("M113").Value = 821.0000199999999
$ws.Range("N113").Value = -5510

